$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data (prices & volume deltas refreshed;
# LEO dropped from the list, Maker added at the bottom, rows 28-51 shifted up).


$ws.Range("D2").Value = "61.759.41"
$ws.Range("E2").Value = "  -0.56%  "

$ws.Range("D3").Value = "3.420.22"
$ws.Range("E3").Value = "  -0.42%  "

$ws.Range("E4").Value = "  +0.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "409.65"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.19"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.633"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.27%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.730"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.139"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "43.54"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000222"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +12.90%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.26"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +6.12%  "

$ws.Range("D14").Value = "3.971.27"
$ws.Range("E14").Value = "  -0.29%  "

$ws.Range("E15").Value = "  +0.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.18"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.59%  "

$ws.Range("D17").Value = "3.432.05"
$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.36"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +7.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.08"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.62%  "

$ws.Range("D20").Value = "61.816.74"
$ws.Range("E20").Value = "  -0.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "491.21"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +32.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "91.75"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.32"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.53"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.33"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +4.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "34.80"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +10.05%  "

$ws.Range("E27").Value = "  +8.78%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.60"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.11%  "

$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "12.13"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.90%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.69"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.85%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.115"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.78%  "

$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.168"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.56%  "

$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "41.75"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.92%  "

$ws.Range("B34").Value = "Dai"
$ws.Range("C34").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "56.16"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +7.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0500"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.65%  "

$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.48"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.63%  "

$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.76"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +18.36%  "

$ws.Range("E40").Value = "  +3.76%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.92"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.65%  "

$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.319"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.82%  "

$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "145.69"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.84%  "

$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.09"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.35"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +8.88%  "

$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.39"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +23.19%  "

$ws.Range("B47").Value = "Celestia"
$ws.Range("C47").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.72"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.61%  "

$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "120.28"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +30.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.03"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.91%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.141"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +16.64%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.135.51"
$ws.Range("E51").Value = "  +1.10%  "
